# Quarterly balance-sheet refresh (Overview sheet)
#
# The source feed rolls its 10-quarter window forward by one quarter:
#   - column D (oldest quarter) is dropped
#   - columns E:M shift left into D:L
#   - column M receives the newly published quarter (Q4 1401/12)
# In addition, the 'read_price' algorithm change recomputed the figures
# for the quarter that now lands in column I (Q4 1400/12), and the publish
# date for the Q3 1401/09 column (now L) was amended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (D8:M8): quarter-period labels
$row8Vals = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value2 = $row8Vals[$i]
}

# Row 9 (D9:M9): publish-date labels (force text so date-like strings
# such as 1399-10-27 aren't reinterpreted as real dates)
$row9Vals = @("'1399-10-27", "'1401-01-31 (8)", "'1400-04-29", "'1400-08-18 (2)", "'1400-10-30", "'1402-01-30 (9)", "'1401-04-30", "'1401-08-02 (2)", "'1401-12-29 (2)", "'1402-01-30 (2)")
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value2 = $row9Vals[$i]
}

# Data rows 12-58 (D:M): rolling quarter shift w/ read_price correction
$rowData = @{}
$rowData[12] = @(433207, 117259, 192026, 618977, 210495, 253515, 957433, 499214, 508520, 435412)
$rowData[13] = @(403950, 368800, 1223360, 388800, 2031448, 0, 0, 0, 0, 0)
$rowData[14] = @(108224, 464357, 1055884, 499389, 504248, 385305, 444563, 303092, 368201, 560921)
$rowData[15] = @(869913, 1142222, 1458559, 1654586, 1897470, 2400230, 0, 2921896, 3183598, 3175528)
$rowData[16] = @(403605, 583408, 745515, 711835, 735697, 270598, 504087, 411182, 482376, 354677)
$rowData[18] = @(2218899, 2676046, 4675344, 3873587, 5379358, 3309648, 1906083, 4135384, 4542695, 4526538)
$rowData[19] = @(668, 1844, 0, 28986, 36298, 13393, 12633, 46730, 46631, 46815)
$rowData[20] = @(874266, 946405, 846405, 1680816, 764325, 3060535, 4113921, 1871113, 2229271, 2891371)
$rowData[22] = @(543378, 579796, 538461, 511117, 472732, 939707, 484858, 1397906, 1400233, 2186586)
$rowData[23] = @(46994, 46995, 47662, 46995, 48994, 46995, 46995, 46995, 46995, 46995)
$rowData[25] = @(0, 0, 0, 0, 0, 0, 2832725, 0, 0, 0)
$rowData[26] = @(1465306, 1575040, 1432528, 2267914, 1322349, 4060630, 7491132, 3362744, 3723130, 5171767)
$rowData[27] = @(3684205, 4251086, 6107872, 6141501, 6701707, 7370278, 9397215, 7498128, 8265825, 9698305)
$rowData[29] = @(362103, 229156, 416653, 642729, 661302, 667745, 997803, 968595, 804787, 678699)
$rowData[31] = @(68002, 163348, 237672, 617312, 251612, 544501, 436149, 212776, 163552, 713511)
$rowData[32] = @(352719, 202252, 362295, 621787, 688000, 394434, 645027, 809151, 659004, 483812)
$rowData[33] = @(56901, 38693, 2338609, 70816, 54948, 57542, 4348542, 688225, 243226, 96131)
$rowData[34] = @(121235, 121235, 121235, 121235, 121235, 0, 121235, 221154, 708987, 358975)
$rowData[35] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[37] = @(960960, 754684, 3476464, 2073879, 1777097, 1664222, 6548756, 2899901, 2579556, 2331128)
$rowData[40] = @(0, 0, 0, 0, 0, 121235, 0, 0, 0, 121235)
$rowData[41] = @(87687, 87343, 125088, 130813, 133727, 128883, 204745, 210492, 213016, 213522)
$rowData[42] = @(87687, 87343, 125088, 130813, 133727, 250118, 204745, 210492, 213016, 334757)
$rowData[43] = @(1048647, 842027, 3601552, 2204692, 1910824, 1914340, 6753501, 3110393, 2792572, 2665885)
$rowData[45] = @(700000, 700000, 700000, 700000, 700000, 700000, 700000, 700000, 1000000, 1000000)
$rowData[47] = @(0, 0, 0, 0, 0, 0, 0, 299496, 0, 0)
$rowData[48] = @(0, -57509, 0, -80960, -80960, -80960, -80960, -80960, -80960, -89021)
$rowData[49] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 32129)
$rowData[50] = @(70000, 70000, 70000, 70000, 70000, 70000, 70000, 70000, 100000, 100000)
$rowData[56] = @(1865558, 2696568, 1736320, 3247769, 4101843, 4766898, 1954674, 3399199, 4454213, 5989312)
$rowData[57] = @(2635558, 3409059, 2506320, 3936809, 4790883, 5455938, 2643714, 4387735, 5473253, 7032420)
$rowData[58] = @(3684205, 4251086, 6107872, 6141501, 6701707, 7370278, 9397215, 7498128, 8265825, 9698305)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value2 = $vals[$i]
    }
}

